$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 8-10 (MuSCs as sending cluster no longer present)
$ws.Range("A8:T10").EntireRow.Delete() | Out-Null

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 66.36304333333334
$ws.Range("H2").Value = 199.08913
$ws.Range("I2").Value = 0.1775372952319303
$ws.Range("J2").Value = 0.1775372952319303
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.866432
$ws.Range("N2").Value = 8.599295999999999
$ws.Range("O2").Value = 0.9456981836489474
$ws.Range("P2").Value = 0.9456981836489475
$ws.Range("Q2").Value = 190.2251510280533
$ws.Range("R2").Value = 1712.02635925248
$ws.Range("S2").Value = 0.1678966976307835
$ws.Range("T2").Value = 0.1678966976307835

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 66.36304333333334
$ws.Range("H3").Value = 199.08913
$ws.Range("I3").Value = 0.1775372952319303
$ws.Range("J3").Value = 0.1775372952319303
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.16459
$ws.Range("N3").Value = 0.49377
$ws.Range("O3").Value = 0.05430181635105255
$ws.Range("P3").Value = 0.05430181635105256
$ws.Range("Q3").Value = 10.92269330223333
$ws.Range("R3").Value = 98.3042397201
$ws.Range("S3").Value = 0.00964059760114688
$ws.Range("T3").Value = 0.00964059760114688

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 296.1091513333333
$ws.Range("H4").Value = 888.327454
$ws.Range("I4").Value = 0.7921640597024409
$ws.Range("J4").Value = 0.7921640597024407
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.866432
$ws.Range("N4").Value = 8.599295999999999
$ws.Range("O4").Value = 0.9456981836489474
$ws.Range("P4").Value = 0.9456981836489475
$ws.Range("Q4").Value = 848.7767468747093
$ws.Range("R4").Value = 7638.990721872383
$ws.Range("S4").Value = 0.7491481124125746
$ws.Range("T4").Value = 0.7491481124125746

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 296.1091513333333
$ws.Range("H5").Value = 888.327454
$ws.Range("I5").Value = 0.7921640597024409
$ws.Range("J5").Value = 0.7921640597024407
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.16459
$ws.Range("N5").Value = 0.49377
$ws.Range("O5").Value = 0.05430181635105255
$ws.Range("P5").Value = 0.05430181635105256
$ws.Range("Q5").Value = 48.73660521795333
$ws.Range("R5").Value = 438.62944696158
$ws.Range("S5").Value = 0.04301594728986617
$ws.Range("T5").Value = 0.04301594728986617

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 11.32556566666667
$ws.Range("H6").Value = 33.976697
$ws.Range("I6").Value = 0.03029864506562886
$ws.Range("J6").Value = 0.03029864506562885
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.866432
$ws.Range("N6").Value = 8.599295999999999
$ws.Range("O6").Value = 0.9456981836489474
$ws.Range("P6").Value = 0.9456981836489475
$ws.Range("Q6").Value = 32.46396384503466
$ws.Range("R6").Value = 292.175674605312
$ws.Range("S6").Value = 0.02865337360558936
$ws.Range("T6").Value = 0.02865337360558935

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 11.32556566666667
$ws.Range("H7").Value = 33.976697
$ws.Range("I7").Value = 0.03029864506562886
$ws.Range("J7").Value = 0.03029864506562885
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.16459
$ws.Range("N7").Value = 0.49377
$ws.Range("O7").Value = 0.05430181635105255
$ws.Range("P7").Value = 0.05430181635105256
$ws.Range("Q7").Value = 1.864074853076667
$ws.Range("R7").Value = 16.77667367769
$ws.Range("S7").Value = 0.001645271460039503
$ws.Range("T7").Value = 0.001645271460039503
